$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Local" results values so F2 and F3 both equal 261178
$ws.Range("F2").Value = 261178
$ws.Range("F3").Value = 261178

# Formulas in F4 (=F2-F3) and F5 (=F4/F2) recalc automatically

# Move the active selection to H22, matching the author's final cursor position
$ws.Range("H22").Select()
